# Adds a new "PAP_DATE_CONFLICT" error row (HEX compilation date parameter)
# to the errors table on the first worksheet, and updates the "Postup" text
# for the existing "SOFTWARE_VERSION_CONFLICT" row to mention the safebytes
# decoding code update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

$updatedPostup = "Skontrolovať, či je štruktúra testovaného záznamu odlišná od známych štruktúr. V prípade, že áno je potrebné od kódu dekódujúceho safebytes túto zmenu zakomponovať."

# Grow the table by one row (gets appended at the end of the table range).
$newRow = $lo.ListRows.Add()

# Remember the old last data row (A18:B18 = 200 / SUCESS) before overwriting it,
# so it can be shifted down to make room for the newly inserted row.
$oldA = $ws.Range("A18").Value2
$oldB = $ws.Range("B18").Value2

$ws.Range("A19").Value = $oldA
$ws.Range("B19").Value = $oldB

# Populate the new row 18 with the PAP_DATE_CONFLICT error details.
$ws.Range("A18").Value = 116
$ws.Range("B18").Value = "PAP_DATE_CONFLICT"
$ws.Range("C18").Value = "Dátum programovania softvéru uložený v safebytes sa nezhoduje s dátumom v hlavičke záznamu."
$ws.Range("D18").Value = $updatedPostup

# Update the "Postup" text of the SOFTWARE_VERSION_CONFLICT row (row 17) to the
# expanded wording that also applies to the new row.
$ws.Range("D17").Value = $updatedPostup

# Match the author's final cursor position.
[void]$ws.Range("C27").Select()
